# Applies the commit "some more cleaning analysis and approach to further
# analysis" to the single 2-column task/progress table in the document.
#
# Net effect (by 1-based table row, column 1 only):
#   Row 3:  "Join the existing datasets"        -> "Clean the names to lower case"
#   Row 4:  "Text mine on the common names"      -> "Join the existing datasets"
#   Row 5:  "Clean the names to lower case"       -> "Text mine on the common names"
#   Row 13: "...use some stats"                   -> "...use some stats vs esrb rating. save as ggplot"
#   Row 18: "Think of an opening question/ theme" -> "Try a gg model thing"
#   Row 19: "Start filling stuff in"               -> "Think of an opening question/ theme"
#   Row 20: "Write  a script of notes for powerpoint" -> "Start filling stuff in"
#   Row 21: "Practise a few times"                 -> "Write  a script of notes for powerpoint"
#   Row 22: "Source pictures"                       -> "Practise a few times"
#   Row 23: "Any surrounding contextual analysis?"  -> "Source pictures"
#   Row 24: "Draw some conclusions"                 -> "Any surrounding contextual analysis?"
#   Row 25: (empty)                                 -> "Draw some conclusions"

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Set-CellPlainText($rowIndex, $text) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $text
}

# --- Simple single-run swaps / shifts -------------------------------------

Set-CellPlainText 3 "Clean the names to lower case"
Set-CellPlainText 4 "Join the existing datasets"
Set-CellPlainText 5 "Text mine on the common names"

Set-CellPlainText 18 "Try a gg model thing"
Set-CellPlainText 19 "Think of an opening question/ theme"
Set-CellPlainText 20 "Start filling stuff in"
Set-CellPlainText 22 "Practise a few times"
Set-CellPlainText 23 "Source pictures"
Set-CellPlainText 24 "Any surrounding contextual analysis?"
Set-CellPlainText 25 "Draw some conclusions"

# --- Row 13: append " vs esrb rating. save as ggplot" (with proofErr tags) -

$row13Cell = $table.Rows.Item(13).Cells.Item(1)
$row13Xml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="7E3468D8" w14:textId="455DDA64" w:rsidR="005E0C23" w:rsidRDefault="005E0C23" w:rsidP="003465F8">' +
    '<w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Check whether I can plot </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ggplot</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> with different geospatial sales volumes vs genres console vs is publisher important use some stats vs </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>esrb</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> rating. save as </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ggplot</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
$row13Cell.Range.InsertXML($row13Xml)

# --- Row 21: becomes the "Write  a script of notes for powerpoint" para ---
# (grammar-checked "Write  a" + spell-checked "powerpoint", moved down from
#  what used to be row 20)

$row21Cell = $table.Rows.Item(21).Cells.Item(1)
$row21Xml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="73173E49" w14:textId="5E475435" w:rsidR="005E0C23" w:rsidRDefault="005E0C23" w:rsidP="003465F8">' +
    '<w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Write  a</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> script of notes for </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>powerpoint</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
$row21Cell.Range.InsertXML($row21Xml)

Write-Output "Applied video games diary task-table edits"
